$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "H2"  = 321
    "I2"  = 899
    "J2"  = 3599
    "K2"  = 19
    "L2"  = 1014
    "M2"  = 63
    "N2"  = 620
    "O2"  = 4
    "P2"  = 15
    "Q2"  = 8
    "R2"  = 65
    "S2"  = 382
    "T2"  = 655
    "U2"  = 44
    "V2"  = 5453
    "W2"  = 2
    "X2"  = 5505
    "Y2"  = 5
    "Z2"  = 88
    "AA2" = 46
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$wb.Save()
